$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 617.1177
$ws.Range("I6").Value = 154.42857
$ws.Range("J6").Value = 941
$ws.Range("K6").Value = 463.28571
$ws.Range("L6").Value = 2823
$ws.Range("M6").Value = -351.28571
$ws.Range("N6").Value = -3047

# Row 43
$ws.Range("H43").Value = 4825.875
$ws.Range("I43").Value = 3753
$ws.Range("J43").Value = 5898.75
$ws.Range("K43").Value = 3753
$ws.Range("L43").Value = 5898.75
$ws.Range("M43").Value = -3684
$ws.Range("N43").Value = -6036.75

# Row 88
$ws.Range("H88").Value = 3941
$ws.Range("I88").Value = 1003
$ws.Range("J88").Value = 4675.5
$ws.Range("K88").Value = 1003
$ws.Range("L88").Value = 4675.5
$ws.Range("M88").Value = -597
$ws.Range("N88").Value = -5487.5

# Row 91
$ws.Range("H91").Value = 3941
$ws.Range("I91").Value = 1003
$ws.Range("J91").Value = 4675.5
$ws.Range("K91").Value = 1003
$ws.Range("L91").Value = 4675.5
$ws.Range("M91").Value = 401
$ws.Range("N91").Value = -7483.5

# Row 138
$ws.Range("H138").Value = 4640
$ws.Range("J138").Value = 4656.3423
$ws.Range("L138").Value = 13969.0269
$ws.Range("N138").Value = -24249.0269

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 45062.203
$ws.Range("I32").Value = 18563.12
$ws.Range("K32").Value = 18563.12
$ws.Range("M32").Value = -18276.12

# Row 133
$ws.Range("H133").Value = 44995
$ws.Range("J133").Value = 44995
$ws.Range("L133").Value = 44995
$ws.Range("N133").Value = -50055

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 55
$ws.Range("H55").Value = 37500
$ws.Range("J55").Value = 37500
$ws.Range("L55").Value = 37500
$ws.Range("N55").Value = -38046

# Row 120
$ws.Range("H120").Value = 33873.5
$ws.Range("J120").Value = 33873.5
$ws.Range("L120").Value = 33873.5
$ws.Range("N120").Value = -43549.5

# Row 134
$ws.Range("H134").Value = 4771.8335
$ws.Range("I134").Value = 4902.8335
$ws.Range("K134").Value = 14708.5005
$ws.Range("M134").Value = -12173.5005

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2609

# Row 65
$ws.Range("H65").Value = 2609

# Row 132
$ws.Range("H132").Value = 23812442
$ws.Range("I132").Value = 3043.2778
$ws.Range("K132").Value = 9129.8334
$ws.Range("M132").Value = -6599.8334

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 320.83334
$ws.Range("I11").Value = 37.5
$ws.Range("J11").Value = 462.5
$ws.Range("K11").Value = 112.5
$ws.Range("L11").Value = 1387.5
$ws.Range("M11").Value = 27.5
$ws.Range("N11").Value = -1667.5

# Row 50
$ws.Range("H50").Value = 1340.7368
$ws.Range("I50").Value = 1687
$ws.Range("J50").Value = 1138.75
$ws.Range("K50").Value = 5061
$ws.Range("L50").Value = 3416.25
$ws.Range("M50").Value = -4580
$ws.Range("N50").Value = -4378.25

# Row 53
$ws.Range("H53").Value = 1340.7368
$ws.Range("I53").Value = 1687
$ws.Range("J53").Value = 1138.75
$ws.Range("K53").Value = 5061
$ws.Range("L53").Value = 3416.25
$ws.Range("M53").Value = -4580
$ws.Range("N53").Value = -4378.25

# Row 68
$ws.Range("H68").Value = 780.3333
$ws.Range("I68").Value = 597
$ws.Range("J68").Value = 911.2857
$ws.Range("K68").Value = 1791
$ws.Range("L68").Value = 2733.8571
$ws.Range("M68").Value = -980
$ws.Range("N68").Value = -4355.8571

# Row 71
$ws.Range("H71").Value = 780.3333
$ws.Range("I71").Value = 597
$ws.Range("J71").Value = 911.2857
$ws.Range("K71").Value = 5373
$ws.Range("L71").Value = 8201.5713
$ws.Range("M71").Value = -1317
$ws.Range("N71").Value = -16313.5713

# Row 82
$ws.Range("H82").Value = 2581.0667
$ws.Range("I82").Value = 1953.3334
$ws.Range("J82").Value = 2999.5557
$ws.Range("K82").Value = 5860.0002
$ws.Range("L82").Value = 8998.667099999999
$ws.Range("M82").Value = -5454.0002
$ws.Range("N82").Value = -9810.667099999999

# Row 85
$ws.Range("H85").Value = 2581.0667
$ws.Range("I85").Value = 1953.3334
$ws.Range("J85").Value = 2999.5557
$ws.Range("K85").Value = 5860.0002
$ws.Range("L85").Value = 8998.667099999999
$ws.Range("M85").Value = -4456.0002
$ws.Range("N85").Value = -11806.6671

# Row 131
$ws.Range("H131").Value = 725131.25
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 746435.1
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 2239305.3
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -2249385.3

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1818.75
$ws.Range("I22").Value = 5050
$ws.Range("J22").Value = 741.6667
$ws.Range("K22").Value = 5050
$ws.Range("L22").Value = 741.6667
$ws.Range("M22").Value = -4755
$ws.Range("N22").Value = -1331.6667

# Row 27
$ws.Range("H27").Value = 1818.75
$ws.Range("I27").Value = 5050
$ws.Range("J27").Value = 741.6667
$ws.Range("K27").Value = 5050
$ws.Range("L27").Value = 741.6667
$ws.Range("M27").Value = -4943
$ws.Range("N27").Value = -955.6667

# Row 46
$ws.Range("H46").Value = 5886.6665
$ws.Range("J46").Value = 4900
$ws.Range("L46").Value = 4900
$ws.Range("N46").Value = -5276

# Row 55
$ws.Range("H55").Value = 977.9231
$ws.Range("I55").Value = 1846.25
$ws.Range("K55").Value = 1846.25
$ws.Range("M55").Value = -1673.25

# Row 68
$ws.Range("H68").Value = 2177.875
$ws.Range("J68").Value = 2343.6
$ws.Range("L68").Value = 2343.6
$ws.Range("N68").Value = -3841.6

# Row 71
$ws.Range("H71").Value = 2177.875
$ws.Range("J71").Value = 2343.6
$ws.Range("L71").Value = 11718
$ws.Range("N71").Value = -19206

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 132
$ws.Range("H132").Value = 3285.1936
$ws.Range("I132").Value = 3447.318
$ws.Range("J132").Value = 2888.889
$ws.Range("K132").Value = 10341.954
$ws.Range("L132").Value = 8666.667000000001
$ws.Range("M132").Value = -7811.954000000002
$ws.Range("N132").Value = -13726.667

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 34019.695
$ws.Range("I132").Value = 4359.8096
$ws.Range("K132").Value = 13079.4288
$ws.Range("M132").Value = -10549.4288
